$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.929.44"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.708.33"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.33"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4023"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +2.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4073"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +0.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.482"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -0.73%  "
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.85"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +1.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08827"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.49"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +6.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.504"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -1.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.039"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.682.48"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.37"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -3.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07177"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.96"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +5.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.273"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.51"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +1.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.922.43"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.340"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.896"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -3.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.427"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +22.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.16"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +1.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.44"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +1.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "143.71"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +5.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.211"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -2.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.267"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +13.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.857.82"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08750"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -1.27%  "
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.03208"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +9.68%  "
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.352"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -2.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.031"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -0.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2878"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +5.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8486"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +7.95%  "
$ws.Range("E40").Value = "  +1.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09464"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +3.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.16"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.478"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.63"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +5.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.729"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +6.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7473"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +3.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.235"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.393"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +5.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.44"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +2.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08418"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +5.59%  "
